$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 165
$ws.Range("A165").Value = "Backend Security Engineer - Hybrid"
$ws.Range("B165").Value = "https://www.dice.com/job-detail/5fb6386a-b996-4f74-befa-eab8ca71dd9b?filters.employmentType=CONTRACTS%7CTHIRD_PARTY&filters.postedDate=ONE&page=2&location=United+States&longitude=-106.5348379&latitude=38.7945952&locationPrecision=Country&countryCode=US&q=Golang"
$ws.Range("C165").Value = "Hybrid in Minneapolis, Minnesota"
$ws.Range("D165").Value = "Contract"
$ws.Range("E165").Value = "Depends on Experience"
$ws.Range("F165").Value = "Hunter Recruiting"

# Row 166
$ws.Range("A166").Value = "Application Support Engineer"
$ws.Range("B166").Value = "https://www.dice.com/job-detail/c20ea71e-eadb-4963-95ca-ea7732f54615?filters.employmentType=CONTRACTS%7CTHIRD_PARTY&filters.postedDate=ONE&page=6&location=United+States&longitude=-106.5348379&latitude=38.7945952&locationPrecision=Country&countryCode=US&q=Golang"
$ws.Range("C166").Value = "Piscataway, New Jersey"
$ws.Range("D166").Value = "Contract"
$ws.Range("E166").Value = "USD 35.00 - 38.00 per hour"
$ws.Range("F166").Value = "Seneca Resources, LLC"
